$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Sheet 1: "Obras en general" ---
$ws1 = $wb.Worksheets.Item("Obras en general")

# Row 8 previously had placeholder empty cells in Q8,T8,U8,V8,W8 (left over from
# a wider selection than actually used). Clear them so no cell record remains.
$ws1.Range("Q8").ClearContents()
$ws1.Range("T8").ClearContents()
$ws1.Range("U8").ClearContents()
$ws1.Range("V8").ClearContents()
$ws1.Range("W8").ClearContents()

# New submission appended as row 9. All values in this sheet are stored as text.
Set-TextValue $ws1.Range("A9") "04/06/2025"
Set-TextValue $ws1.Range("B9") "Ingeniero"
Set-TextValue $ws1.Range("C9") "Físico"
Set-TextValue $ws1.Range("D9") "2"
Set-TextValue $ws1.Range("E9") "Obra nueva"
Set-TextValue $ws1.Range("F9") "JUANI GALLO"
Set-TextValue $ws1.Range("G9") "FIRU LUQUE"
Set-TextValue $ws1.Range("H9") "FRANCISCO DE HARO 2745"
Set-TextValue $ws1.Range("I9") "25817/G/2025"
Set-TextValue $ws1.Range("K9") "151818"
Set-TextValue $ws1.Range("L9") "15000"
Set-TextValue $ws1.Range("N9") "7000"
Set-TextValue $ws1.Range("O9") "4000"
Set-TextValue $ws1.Range("R9") "No pagado"
Set-TextValue $ws1.Range("S9") "No pagado"
Set-TextValue $ws1.Range("Y9") "3764251817"

# --- Sheet 2: "Informes técnicos" ---
$ws2 = $wb.Worksheets.Item("Informes técnicos")

Set-TextValue $ws2.Range("A2") "04/06/2025"
Set-TextValue $ws2.Range("B2") "Licenciado"
Set-TextValue $ws2.Range("C2") "Físico"
Set-TextValue $ws2.Range("D2") "7"
Set-TextValue $ws2.Range("E2") "Plan de Contingencia"
Set-TextValue $ws2.Range("F2") "EVACUACION"
Set-TextValue $ws2.Range("G2") "LAUDIN JORGE"
Set-TextValue $ws2.Range("H2") "ARMOA ESTELA"
Set-TextValue $ws2.Range("I2") "8000"
Set-TextValue $ws2.Range("J2") "No pagado"

# These columns are submitted blank on the form (no value to record).
$ws2.Range("K2").Value = ""
$ws2.Range("L2").Value = ""
$ws2.Range("M2").Value = ""
$ws2.Range("N2").Value = ""

Set-TextValue $ws2.Range("P2") "3764251817"
